$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price-log entry was added for "Terminal La Palmera de La
# Serena - Uva": insert a brand-new row at row 65 (this pushes the
# existing rows 65-132 down to 66-133, matching the new dimension
# A1:T133) and populate it with the new record's data.
$ws.Rows(65).Insert()

$ws.Range("A65").Value = 8
$ws.Range("B65").Value = "Terminal La Palmera de La Serena"
$ws.Range("C65").Value = "Coquimbo"
$ws.Range("D65").Value = 44944
$ws.Range("E65").Value = 4
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100109
$ws.Range("H65").Value = "Uva"
$ws.Range("I65").Value = 100109001
$ws.Range("J65").Value = "Uva"
$ws.Range("K65").Value = "Flame Seedless"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 1600
$ws.Range("N65").Value = 8000
$ws.Range("O65").Value = 9000
$ws.Range("P65").Value = 8500
$ws.Range("Q65").Value = "`$/bandeja 10 kilos"
$ws.Range("R65").Value = "Provincia del Elquí"
$ws.Range("S65").Value = 850
$ws.Range("T65").Value = 10
